# kelp_reproductive_timing.xlsx — Nov 2023 data + "time_search" reproduction sheet
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Sheet1: append the Nov/Oct? new data row (A20:C20) -> 2023-10-20 observation
# ---------------------------------------------------------------------------
$ws1.Range("A20").Value = 2023
$ws1.Range("B20").Value = 10
$ws1.Range("C20").Value = 20

# ---------------------------------------------------------------------------
# 2) Add "Sheet2" right after "Sheet1" — a reproduction of the dataset that
#    adds a type_of_data column and a new Nov-17-2023 "time_search" row.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)

# Column widths (tuned to match the source workbook's layout)
$ws2.Columns("B").ColumnWidth = 11.75006
$ws2.Columns("C").ColumnWidth = 12.58349
$ws2.Columns("D").ColumnWidth = 8.58349
$ws2.Columns("E").ColumnWidth = 9.75006
$ws2.Columns("F").ColumnWidth = 5.41683
$ws2.Range("G1:H1").ColumnWidth = 7.0834

# Header row
$ws2.Range("A1").Value = "year"
$ws2.Range("B1").Value = "month"
$ws2.Range("C1").Value = "day"
$ws2.Range("D1").Value = "saccharina_latissima"
$ws2.Range("E1").Value = "nereocystis_luetkeana"
$ws2.Range("F1").Value = "alaria_marginata"
$ws2.Range("G1").Value = "costaria_costata"
$ws2.Range("H1").Value = "type_of_data"

# Data rows 2-20: same observations as Sheet1, each tagged presence_abscence
$rows = @(
    @(2022, 4,  20, 1,    $null),
    @(2022, 5,  19, $null, 1),
    @(2022, 6,  15, $null, $null),
    @(2022, 7,  14, $null, 1),
    @(2022, 8,  $null, $null, $null),
    @(2022, 9,  23, $null, 1),
    @(2022, 10, 13, 1,    $null),
    @(2022, 10, 27, 1,    $null),
    @(2022, 12, 25, 1,    $null),
    @(2023, 1,  $null, $null, $null),
    @(2023, 2,  18, 1,    $null),
    @(2023, 3,  $null, $null, $null),
    @(2023, 4,  $null, $null, $null),
    @(2023, 5,  8,  1,    1),
    @(2023, 6,  6,  1,    1),
    @(2023, 7,  4,  $null, 1),
    @(2023, 8,  1,  1,    1),
    @(2023, 9,  1,  1,    $null),
    @(2023, 10, 20, $null, $null)
)

$r = 2
foreach ($row in $rows) {
    $ws2.Range("A$r").Value = $row[0]
    $ws2.Range("B$r").Value = $row[1]
    if ($row[2] -ne $null) { $ws2.Range("C$r").Value = $row[2] }
    if ($row[3] -ne $null) { $ws2.Range("D$r").Value = $row[3] }
    if ($row[4] -ne $null) { $ws2.Range("E$r").Value = $row[4] }
    $ws2.Range("H$r").Value = "presence_abscence"
    $r++
}

# Row 21: new Nov 17 2023 "time_search" observation, highlighted yellow
$ws2.Range("A21").Value = 2023
$ws2.Range("B21").Value = 11
$ws2.Range("C21").Value = 17
$ws2.Range("D21").Value = 203
$ws2.Range("E21").Value = 0
$ws2.Range("F21").Value = 0
$ws2.Range("G21").Value = 0
$ws2.Range("H21").Value = "time_search"
$ws2.Range("A21:H21").Interior.Color = 65535

# ---------------------------------------------------------------------------
# 3) View state: freeze/scroll Sheet1 back to the top, select all of it, then
#    leave Sheet2 active with its own selection — matching the saved file.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 2
$ws1.Cells.Select()

$ws2.Range("M9").Select()
$ws2.Activate()
